# Updated cryptos list (price & 1h volume % change) across rows 2-51
# Source: automated crypto price refresh (GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values are plain numeric strings (e.g. "251.94")
# would otherwise be auto-converted from text to a Number by Excel's normal
# type inference when written via .Value. Force those specific cells to a
# Text number format before the write so the value is stored as a string,
# matching the source data (prices are free-form display strings, not real
# numbers - see e.g. "30.353.03" elsewhere in the same column).
$numericLookingPriceCells = @("D5","D6","D8","D9","D10","D11","D12","D14","D15","D16","D18","D19","D20","D22","D23","D24","D25","D26","D27","D28","D29","D30","D32","D33","D34","D35","D36","D38","D39","D40","D41","D42","D43","D44","D46","D47","D48","D49","D50","D51")
foreach ($addr in $numericLookingPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply every updated Price (D) / Volume(1h) (E) cell value
$ws.Range("D2").Value = '30.353.03'
$ws.Range("E2").Value = '  -2.77%  '
$ws.Range("D3").Value = '1.940.58'
$ws.Range("E3").Value = '  -2.91%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '251.94'
$ws.Range("E5").Value = '  -1.88%  '
$ws.Range("D6").Value = '0.7184'
$ws.Range("E6").Value = '  -7.57%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '0.3344'
$ws.Range("E8").Value = '  -3.58%  '
$ws.Range("D9").Value = '28.58'
$ws.Range("E9").Value = '  -0.49%  '
$ws.Range("D10").Value = '0.07313'
$ws.Range("E10").Value = '  +4.38%  '
$ws.Range("D11").Value = '0.8155'
$ws.Range("E11").Value = '  -4.50%  '
$ws.Range("D12").Value = '0.08149'
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("D13").Value = '1.939.18'
$ws.Range("E13").Value = '  -2.91%  '
$ws.Range("D14").Value = '5.498'
$ws.Range("E14").Value = '  -1.82%  '
$ws.Range("D15").Value = '95.30'
$ws.Range("E15").Value = '  -5.42%  '
$ws.Range("D16").Value = '14.89'
$ws.Range("E16").Value = '  -4.34%  '
$ws.Range("D17").Value = '30.367.61'
$ws.Range("E17").Value = '  -2.75%  '
$ws.Range("D18").Value = '0.000008330'
$ws.Range("E18").Value = '  +5.30%  '
$ws.Range("D19").Value = '253.54'
$ws.Range("E19").Value = '  -7.52%  '
$ws.Range("D20").Value = '5.864'
$ws.Range("E20").Value = '  -1.58%  '
$ws.Range("D21").Value = '2.194.25'
$ws.Range("E21").Value = '  -2.92%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '6.974'
$ws.Range("E24").Value = '  -1.93%  '
$ws.Range("D25").Value = '9.839'
$ws.Range("E25").Value = '  -1.93%  '
$ws.Range("D26").Value = '161.19'
$ws.Range("E26").Value = '  -2.17%  '
$ws.Range("D27").Value = '2.404'
$ws.Range("E27").Value = '  +3.37%  '
$ws.Range("D28").Value = '19.40'
$ws.Range("E28").Value = '  -2.63%  '
$ws.Range("D29").Value = '0.1311'
$ws.Range("E29").Value = '  -10.74%  '
$ws.Range("D30").Value = '1.579'
$ws.Range("E30").Value = '  -1.73%  '
$ws.Range("E31").Value = '  -0.92%  '
$ws.Range("D32").Value = '4.490'
$ws.Range("E32").Value = '  -2.66%  '
$ws.Range("D33").Value = '4.258'
$ws.Range("E33").Value = '  -3.85%  '
$ws.Range("D34").Value = '0.05283'
$ws.Range("E34").Value = '  +1.35%  '
$ws.Range("D35").Value = '1.270'
$ws.Range("E35").Value = '  +3.35%  '
$ws.Range("D36").Value = '0.7644'
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = '0.01997'
$ws.Range("E38").Value = '  -0.17%  '
$ws.Range("D39").Value = '2.845'
$ws.Range("E39").Value = '  -2.62%  '
$ws.Range("D40").Value = '80.87'
$ws.Range("E40").Value = '  +1.77%  '
$ws.Range("D41").Value = '6.565'
$ws.Range("E41").Value = '  -2.18%  '
$ws.Range("D42").Value = '0.4576'
$ws.Range("E42").Value = '  -2.27%  '
$ws.Range("D43").Value = '2.030'
$ws.Range("E43").Value = '  -5.37%  '
$ws.Range("D44").Value = '0.8483'
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").Value = '103.28'
$ws.Range("E46").Value = '  -2.62%  '
$ws.Range("D47").Value = '9.859'
$ws.Range("E47").Value = '  -0.43%  '
$ws.Range("D48").Value = '7.397'
$ws.Range("E48").Value = '  -4.18%  '
$ws.Range("D49").Value = '37.18'
$ws.Range("E49").Value = '  +1.00%  '
$ws.Range("D50").Value = '0.4184'
$ws.Range("E50").Value = '  -2.86%  '
$ws.Range("D51").Value = '1.503'
$ws.Range("E51").Value = '  -1.53%  '

# Restore the default cell style on the coerced cells so only their value
# (not their formatting) differs from the original workbook.
foreach ($addr in $numericLookingPriceCells) {
    $ws.Range($addr).Style = "Normal"
}
